# Scheduled runner: refresh live market-price-derived columns (H-N)
# on the Brynhildr_Profits workbook. Values come from the latest
# Universalis price snapshot; only numeric cells are rewritten.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 433.66666
$ws.Range("I39").Value = 28.25
$ws.Range("K39").Value = 84.75
$ws.Range("M39").Value = 211.25
$ws.Range("H62").Value = 4080.4443
$ws.Range("I62").Value = 3165.889
$ws.Range("K62").Value = 3165.889
$ws.Range("M62").Value = -2541.889
$ws.Range("H64").Value = 5623.3335
$ws.Range("I64").Value = 4435
$ws.Range("K64").Value = 4435
$ws.Range("M64").Value = -4187
$ws.Range("H65").Value = 4080.4443
$ws.Range("I65").Value = 3165.889
$ws.Range("K65").Value = 15829.445
$ws.Range("M65").Value = -12709.445
$ws.Range("H67").Value = 5623.3335
$ws.Range("I67").Value = 4435
$ws.Range("K67").Value = 4435
$ws.Range("M67").Value = -3577
$ws.Range("H98").Value = 2017.0605
$ws.Range("I98").Value = 1962.2069
$ws.Range("K98").Value = 1962.2069
$ws.Range("M98").Value = -464.2068999999999
$ws.Range("H122").Value = 2017.0605
$ws.Range("I122").Value = 1962.2069
$ws.Range("K122").Value = 5886.620699999999
$ws.Range("M122").Value = -3436.620699999999
$ws.Range("H132").Value = 11299.368
$ws.Range("I132").Value = 11871.833
$ws.Range("K132").Value = 35615.499
$ws.Range("M132").Value = -33085.499
$ws.Range("H138").Value = 4768.4644
$ws.Range("I138").Value = 5646.3335
$ws.Range("J138").Value = 4352.6313
$ws.Range("K138").Value = 16939.0005
$ws.Range("L138").Value = 13057.8939
$ws.Range("M138").Value = -11799.0005
$ws.Range("N138").Value = -23337.8939
$ws.Range("H141").Value = 12476.808
$ws.Range("I141").Value = 4600
$ws.Range("J141").Value = 13908.954
$ws.Range("K141").Value = 13800
$ws.Range("L141").Value = 41726.862
$ws.Range("M141").Value = -8620
$ws.Range("N141").Value = -52086.862

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1544
$ws.Range("I2").Value = 1280.4
$ws.Range("K2").Value = 1280.4
$ws.Range("M2").Value = -1167.4
$ws.Range("H61").Value = 6669254.5
$ws.Range("I61").Value = 2678.4614
$ws.Range("K61").Value = 2678.4614
$ws.Range("M61").Value = -2466.4614
$ws.Range("H116").Value = 1544
$ws.Range("I116").Value = 1280.4
$ws.Range("K116").Value = 1280.4
$ws.Range("M116").Value = 1013.6
$ws.Range("H132").Value = 2622.15
$ws.Range("I132").Value = 1065.05
$ws.Range("J132").Value = 5736.35
$ws.Range("K132").Value = 3195.15
$ws.Range("L132").Value = 17209.05
$ws.Range("M132").Value = -665.1499999999996
$ws.Range("N132").Value = -22269.05
$ws.Range("H136").Value = 6669254.5
$ws.Range("I136").Value = 2678.4614
$ws.Range("K136").Value = 8035.3842
$ws.Range("M136").Value = -5485.3842
$ws.Range("H139").Value = 100715
$ws.Range("J139").Value = 100715
$ws.Range("L139").Value = 100715
$ws.Range("N139").Value = -110995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1544
$ws.Range("I3").Value = 1280.4
$ws.Range("K3").Value = 1280.4
$ws.Range("M3").Value = -1166.4
$ws.Range("H134").Value = 2033658
$ws.Range("I134").Value = 1191.2632
$ws.Range("J134").Value = 27778236
$ws.Range("K134").Value = 3573.7896
$ws.Range("L134").Value = 83334708
$ws.Range("M134").Value = -1038.7896
$ws.Range("N134").Value = -83339778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5090670
$ws.Range("I31").Value = 2019126
$ws.Range("K31").Value = 2019126
$ws.Range("M31").Value = -2018831
$ws.Range("H34").Value = 5090670
$ws.Range("I34").Value = 2019126
$ws.Range("K34").Value = 2019126
$ws.Range("M34").Value = -2018924
$ws.Range("H58").Value = 6634583.5
$ws.Range("I58").Value = 8773841
$ws.Range("J58").Value = 3924857.8
$ws.Range("K58").Value = 8773841
$ws.Range("L58").Value = 3924857.8
$ws.Range("M58").Value = -8773638
$ws.Range("N58").Value = -3925263.8
$ws.Range("H99").Value = 30815
$ws.Range("I99").Value = 41881
$ws.Range("K99").Value = 41881
$ws.Range("M99").Value = -40383
$ws.Range("H126").Value = 30815
$ws.Range("I126").Value = 41881
$ws.Range("K126").Value = 125643
$ws.Range("M126").Value = -123173
$ws.Range("H132").Value = 3141.6667
$ws.Range("I132").Value = 2914.2144
$ws.Range("J132").Value = 3937.75
$ws.Range("K132").Value = 8742.643199999999
$ws.Range("L132").Value = 11813.25
$ws.Range("M132").Value = -6212.643199999999
$ws.Range("N132").Value = -16873.25
$ws.Range("H134").Value = 3845.8115
$ws.Range("I134").Value = 2765.9033
$ws.Range("J134").Value = 4726.7896
$ws.Range("K134").Value = 8297.7099
$ws.Range("L134").Value = 14180.3688
$ws.Range("M134").Value = -5762.7099
$ws.Range("N134").Value = -19250.3688
$ws.Range("H136").Value = 6634583.5
$ws.Range("I136").Value = 8773841
$ws.Range("J136").Value = 3924857.8
$ws.Range("K136").Value = 26321523
$ws.Range("L136").Value = 11774573.4
$ws.Range("M136").Value = -26318973
$ws.Range("N136").Value = -11779673.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 427
$ws.Range("I7").Value = 574.4
$ws.Range("J7").Value = 279.6
$ws.Range("K7").Value = 1723.2
$ws.Range("L7").Value = 838.8000000000001
$ws.Range("M7").Value = -1611.2
$ws.Range("N7").Value = -1062.8
$ws.Range("H51").Value = 133334330
$ws.Range("I51").Value = 133334330
$ws.Range("K51").Value = 400002990
$ws.Range("M51").Value = -400002530
$ws.Range("H131").Value = 4230.8037
$ws.Range("I131").Value = 557.2143
$ws.Range("J131").Value = 5620.811
$ws.Range("K131").Value = 1671.6429
$ws.Range("L131").Value = 16862.433
$ws.Range("M131").Value = 3368.3571
$ws.Range("N131").Value = -26942.433
$ws.Range("H133").Value = 8256.25
$ws.Range("I133").Value = 4828.778
$ws.Range("J133").Value = 11060.546
$ws.Range("K133").Value = 14486.334
$ws.Range("L133").Value = 33181.638
$ws.Range("M133").Value = -9426.334000000001
$ws.Range("N133").Value = -43301.638
$ws.Range("H140").Value = 3004.3076
$ws.Range("I140").Value = 2570.087
$ws.Range("K140").Value = 7710.261
$ws.Range("M140").Value = -2530.261

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 9007.5
$ws.Range("I55").Value = 3015
$ws.Range("K55").Value = 3015
$ws.Range("M55").Value = -2688
$ws.Range("H132").Value = 11085.595
$ws.Range("I132").Value = 6084.353
$ws.Range("J132").Value = 67766.336
$ws.Range("K132").Value = 18253.059
$ws.Range("L132").Value = 203299.008
$ws.Range("M132").Value = -15723.059
$ws.Range("N132").Value = -208359.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 38499.668
$ws.Range("I36").Value = 49499
$ws.Range("J36").Value = 33000
$ws.Range("K36").Value = 49499
$ws.Range("L36").Value = 33000
$ws.Range("M36").Value = -48937
$ws.Range("N36").Value = -34124
$ws.Range("H132").Value = 904078.8
$ws.Range("I132").Value = 1518179
$ws.Range("K132").Value = 4554537
$ws.Range("M132").Value = -4552007

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 2499.5
$ws.Range("J13").Value = 2499.5
$ws.Range("L13").Value = 2499.5
$ws.Range("N13").Value = -2779.5
$ws.Range("H113").Value = 1378.9487
$ws.Range("J113").Value = 1670.2632
$ws.Range("L113").Value = 5010.7896
$ws.Range("N113").Value = -9350.7896
$ws.Range("H136").Value = 5512278
$ws.Range("I136").Value = 2808297.8
$ws.Range("K136").Value = 8424893.399999999
$ws.Range("M136").Value = -8422343.399999999
